$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo: B10 should be "aa" instead of "a"
$ws.Range("B10").Value = "aa"

# Update selection to reflect active cell B10
$ws.Range("B10").Select()
